# Crops energy content.xlsx - add a curated "Basic database" sheet in front of
# the original data sheet (renamed "Expanded database"), and tweak a handful of
# formatting details on the original sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the original (only) sheet and insert a new sheet before it.
# ---------------------------------------------------------------------------
$expanded = $wb.Worksheets.Item(1)
$expanded.Name = "Expanded database"

$basic = $wb.Worksheets.Add()
$basic.Name = "Basic database"

# ---------------------------------------------------------------------------
# 2. Populate "Basic database" - header row + a curated subset of 9 food
#    items taken from the expanded list, plus a brand-new "Microgreen mix AH"
#    entry.
# ---------------------------------------------------------------------------
$data = New-Object 'object[,]' 10,12

# Header
$data[0,0]="NEVO-code";  $data[0,1]="Voedingsmiddel"; $data[0,2]="Food item"
$data[0,3]="Energie (kJ)"; $data[0,4]="Energie (kcal)"; $data[0,5]="Eiwit totaal (g)"
$data[0,6]="Koolhydraten totaal (g)"; $data[0,7]="Vet totaal (g)"; $data[0,8]="Voedingsvezel totaal (g)"
$data[0,9]="Alcohol totaal (g)"; $data[0,10]="Water (g)"; $data[0,11]="As (g)"

# Lettuce
$data[1,0]=1399; $data[1,1]="Sla ijsberg- rauw"; $data[1,2]="Lettuce"
$data[1,3]=62; $data[1,4]=15; $data[1,5]=0.9; $data[1,6]=1.7; $data[1,7]=0.3
$data[1,8]=1.1; $data[1,9]=0; $data[1,10]=95.7; $data[1,11]=0.5

# Endive
$data[2,0]=7; $data[2,1]="Andijvie rauw"; $data[2,2]="Endive "
$data[2,3]=71; $data[2,4]=17; $data[2,5]=1.3; $data[2,6]=1.2; $data[2,7]=0.4
$data[2,8]=1.8; $data[2,9]=0; $data[2,10]=95.3; $data[2,11]=0.6

# Spinach
$data[3,0]=51; $data[3,1]="Spinazie rauw"; $data[3,2]="Spinach"
$data[3,3]=108; $data[3,4]=26; $data[3,5]=3.2; $data[3,6]=0.9; $data[3,7]=0.6
$data[3,8]=2; $data[3,9]=0; $data[3,10]=91.9; $data[3,11]=1.4

# Bean sprouts
$data[4,0]=58; $data[4,1]="Tauge rauw"; $data[4,2]="Bean sprouts"
$data[4,3]=97; $data[4,4]=23; $data[4,5]=2.2; $data[4,6]=2.7; $data[4,7]=0.1
$data[4,8]=1.2; $data[4,9]=0; $data[4,10]=93; $data[4,11]=0.3

# Parsley fresh
$data[5,0]=128; $data[5,1]="Peterselie vers"; $data[5,2]="Parsley fresh"
$data[5,3]=125; $data[5,4]=30; $data[5,5]=4; $data[5,6]=1; $data[5,7]=0
$data[5,8]=5; $data[5,9]=0; $data[5,10]=87; $data[5,11]="--"

# Kale curly
$data[6,0]=959; $data[6,1]="Kool boeren- rauw"; $data[6,2]="Kale curly"
$data[6,3]=193; $data[6,4]=46; $data[6,5]=4; $data[6,6]=4; $data[6,7]=1
$data[6,8]=2.5; $data[6,9]=0; $data[6,10]=86; $data[6,11]="--"

# Basil fresh
$data[7,0]=2177; $data[7,1]="Basilicum vers"; $data[7,2]="Basil fresh"
$data[7,3]=200; $data[7,4]=48; $data[7,5]=3.1; $data[7,6]=5.1; $data[7,7]=0.8
$data[7,8]=3.9; $data[7,9]=0; $data[7,10]=$null; $data[7,11]=$null

# Rucola
$data[8,0]=2736; $data[8,1]="Sla rucola rauw"; $data[8,2]="Rucola"
$data[8,3]=98; $data[8,4]=23; $data[8,5]=3.5; $data[8,6]=0.1; $data[8,7]=0.5
$data[8,8]=2.2; $data[8,9]=0; $data[8,10]=87; $data[8,11]="--"

# Microgreen mix AH (new entry, only partial data like the others' "--" rows)
$data[9,0]=$null; $data[9,1]=$null; $data[9,2]="Microgreen mix AH"
$data[9,3]=$null; $data[9,4]=32; $data[9,5]=$null; $data[9,6]=$null; $data[9,7]=$null
$data[9,8]=$null; $data[9,9]=$null; $data[9,10]=96.7; $data[9,11]=0.3

$basic.Range("A1:L10").Value2 = $data

# ---------------------------------------------------------------------------
# 3. Basic formatting pass on "Basic database" - column widths similar to
#    the source workbook, default row height, and the active selection.
# ---------------------------------------------------------------------------
$basic.Range("A1:L10").Font.Name = "Arial"
$basic.Range("A1:L10").Font.Size = 10

$basic.Columns.Item(1).ColumnWidth = 10.86
$basic.Columns.Item(2).ColumnWidth = 23
$basic.Columns.Item(3).ColumnWidth = 17
$basic.Columns.Item(4).ColumnWidth = 10.86
$basic.Columns.Item(5).ColumnWidth = 12.43
$basic.Columns.Item(6).ColumnWidth = 12.86
$basic.Columns.Item(7).ColumnWidth = 20
$basic.Columns.Item(8).ColumnWidth = 11.71
$basic.Columns.Item(9).ColumnWidth = 21
$basic.Columns.Item(10).ColumnWidth = 9.43

$basic.Range("D23").Select()

# ---------------------------------------------------------------------------
# 4. Formatting tweaks on "Expanded database": column C styling, row heights,
#    and the saved selection / scroll position.
# ---------------------------------------------------------------------------
$expanded.Columns.Item(3).Font.Name = "Arial"
$expanded.Columns.Item(3).Font.Size = 10

for ($r = 1; $r -le 25; $r++) {
    if ($r -ne 22) {
        $expanded.Rows.Item($r).RowHeight = 12.75
    }
}

$expanded.Range("A34:M44").Select()
